# bug fixes for import and change text restaurant intoxication to cancelled

$wb = $excel.ActiveWorkbook
$items = $wb.Worksheets.Item("Items")
$variations = $wb.Worksheets.Item("Variations")

# --- Items sheet: header row text fix ---
$items.Range("B1").Value = "Name"
$items.Range("I1").Value = "description"

# --- Items sheet: row 2 updates ---
$items.Range("B2").Value = "New Drink 2"
$items.Range("I2").Value = "hello description 1"
$items.Range("J2").Value = 1

# --- Items sheet: row 3 updates ---
$items.Range("C3").Value = "Spiritss"
$items.Range("I3").Value = "hello description 2"
$items.Range("J3").Value = 1

# --- Items sheet: add new rows 4-6, carrying row 3's per-cell formatting ---
$items.Range("A3:L3").Copy()
$items.Range("A4:L4").PasteSpecial(-4122)
$items.Range("A3:L3").Copy()
$items.Range("A5:L5").PasteSpecial(-4122)
$items.Range("A3:L3").Copy()
$items.Range("A6:L6").PasteSpecial(-4122)

$items.Range("A4").Value = 3
$items.Range("B4").Value = "New Drink 2"
$items.Range("C4").Value = "Spiritss"
$items.Range("D4").Value = 110
$items.Range("E4").Value = "test2"
$items.Range("F4").Value = "Aus"
$items.Range("G4").Value = 1997
$items.Range("H4").Value = "Wine"
$items.Range("I4").Value = "hello description 3"
$items.Range("J4").Value = 0
$items.Range("K4").Value = 1
$items.Range("L4").Value = 0

$items.Range("A5").Value = 4
$items.Range("B5").Value = "New Drink 2"
$items.Range("C5").Value = "Spiritss"
$items.Range("D5").Value = 110
$items.Range("E5").Value = "test2"
$items.Range("F5").Value = "Aus"
$items.Range("G5").Value = 1997
$items.Range("H5").Value = "Wine"
$items.Range("I5").Value = "hello description 4"
$items.Range("J5").Value = 0
$items.Range("K5").Value = 1
$items.Range("L5").Value = 0

$items.Range("A6").Value = 5
$items.Range("B6").Value = "New Drink 2"
$items.Range("C6").Value = "Spiritss"
$items.Range("D6").Value = 110
$items.Range("E6").Value = "test2"
$items.Range("F6").Value = "Aus"
$items.Range("G6").Value = 1997
$items.Range("H6").Value = "Wine"
$items.Range("I6").Value = "hello description 5"
$items.Range("J6").Value = 0
$items.Range("K6").Value = 1
$items.Range("L6").Value = 0

# column I (description) got a bit wider to fit the new text
$items.Columns.Item(9).ColumnWidth = 16.85

# header row is shorter now (less wrapped text)
$items.Rows.Item(1).RowHeight = 90

# selection/active-view tweaks: focus moves to the Items sheet
$variations.Range("G5").Select()
$items.Range("I13").Select()
